# Update English (en_US, column C) lines of the Kazdel/Hoederer/Ines dialogue
# to replace curly "smart" double quotes with straight single quotes,
# per commit "update on 20210731 画中人".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value  = "[name=`"Hoederer`"]  How many 'nobles' do you think there are in Kazdel? How many of them are Sarkaz who gave themselves some stupid title in an even stupider war?`n"
$ws.Range("C65").Value  = "[name=`"Hoederer`"]  What do you think about Her Majesty? And Kal'tsit? And especially... that 'Doctor?'`n"
$ws.Range("C72").Value  = "[name=`"Ines`"]  But, I might’ve accidentally taken a peek at Her Majesty. My Arts can’t look directly into peoples’ minds. More like I can get a feel for them by looking over their 'shadows'...`n"
$ws.Range("C78").Value  = "[name=`"Ines`"]  She sees us all as equals. She didn’t call us 'devils.'`n"
$ws.Range("C80").Value  = "[name=`"Hoederer`"]  What about the one they just call 'The Doctor?'`n"
$ws.Range("C101").Value = "[name=`"Hoederer`"]  Or do you think the Regent, Her Majesty’s brother, will open his heart, say 'all is forgiven,' and keep letting us raid the Lateranos?`n"
$ws.Range("C121").Value = "[name=`"W`"]  Yeah, the 'alive' part is trickier than you’d think. Guy offed himself before I could even murk him. What’s a girl supposed to do?`n"
$ws.Range("C149").Value = "Uhh, probably better than you do? I’m the one actually down on the battlefield putting in work, 'Dr. K'~`n"
$ws.Range("C169").Value = "She and I are the same kind of person. If she can even be considered a 'person,' that is. That’s one of a whole bunch of mysteries floating around her.`n"
